$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semilla 11")

$ws.Range("B9").Value = "708603165"
$ws.Range("C9").Value = "3046010569"

$ws.Range("B10").Value = "325869013"
$ws.Range("C10").Value = "3052749177"

$ws.Range("B11").Value = "702923575"
$ws.Range("C11").Value = "3046010569"

$ws.Range("C12").Value = "3052754285"

$ws.Range("C13").Value = "3052754289"

$ws.Range("C14").Value = "3052754293"

$ws.Range("C15").Value = "3052754321"

$ws.Range("D18").Select()
